$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157, shifting existing rows 157..241 down to 158..242
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with its values
$ws.Cells.Item(157, 1).Value = 4
$ws.Cells.Item(157, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(157, 3).Value = "Los Lagos"
$ws.Cells.Item(157, 4).Value = 44460
$ws.Cells.Item(157, 5).Value = 10
$ws.Cells.Item(157, 6).Value = "Fruta"
$ws.Cells.Item(157, 7).Value = 100108
$ws.Cells.Item(157, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(157, 9).Value = 100108006
$ws.Cells.Item(157, 10).Value = "Plátano"
$ws.Cells.Item(157, 11).Value = "Sin especificar"
$ws.Cells.Item(157, 12).Value = "Primera Pintón"
$ws.Cells.Item(157, 13).Value = 400
$ws.Cells.Item(157, 14).Value = 18500
$ws.Cells.Item(157, 15).Value = 18500
$ws.Cells.Item(157, 16).Value = 18500
$ws.Cells.Item(157, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(157, 18).Value = "Ecuador"
$ws.Cells.Item(157, 19).Value = 925
$ws.Cells.Item(157, 20).Value = 20
